$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.047761014565285
$ws.Cells.Item(2, 4).Value = 1.04858369611772
$ws.Cells.Item(2, 5).Value = 1.061476079636409
$ws.Cells.Item(2, 6).Value = 1.069092229906419
$ws.Cells.Item(2, 9).Value = 1.047367400325523
$ws.Cells.Item(2, 10).Value = 1.052808364554159
$ws.Cells.Item(2, 11).Value = 1.051343135170294
$ws.Cells.Item(2, 12).Value = 1.064200037213684
$ws.Cells.Item(2, 13).Value = 1.071795661885899
$ws.Cells.Item(2, 14).Value = 1.054303473190633
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.048713810897527
$ws.Cells.Item(3, 4).Value = 1.049298136647545
$ws.Cells.Item(3, 5).Value = 1.062420903950883
$ws.Cells.Item(3, 6).Value = 1.070119402645932
$ws.Cells.Item(3, 9).Value = 1.047653923427502
$ws.Cells.Item(3, 10).Value = 1.053409748226745
$ws.Cells.Item(3, 11).Value = 1.051869902944492
$ws.Cells.Item(3, 12).Value = 1.064959161825411
$ws.Cells.Item(3, 13).Value = 1.07263841461085
$ws.Cells.Item(3, 14).Value = 1.054905710897013
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.049330711910519
$ws.Cells.Item(4, 4).Value = 1.049760670503712
$ws.Cells.Item(4, 5).Value = 1.063033009318421
$ws.Cells.Item(4, 6).Value = 1.070784931949798
$ws.Cells.Item(4, 9).Value = 1.047838209235339
$ws.Cells.Item(4, 10).Value = 1.05379862556187
$ws.Cells.Item(4, 11).Value = 1.052210326822326
$ws.Cells.Item(4, 12).Value = 1.065450479176373
$ws.Cells.Item(4, 13).Value = 1.07318399904279
$ws.Cells.Item(4, 14).Value = 1.055295140482559
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.049590146791748
$ws.Cells.Item(5, 4).Value = 1.049955176735991
$ws.Cells.Item(5, 5).Value = 1.063290514669869
$ws.Cells.Item(5, 6).Value = 1.071064930035245
$ws.Cells.Item(5, 9).Value = 1.047915415833096
$ws.Cells.Item(5, 10).Value = 1.053962047000609
$ws.Cells.Item(5, 11).Value = 1.052353337240127
$ws.Cells.Item(5, 12).Value = 1.065657055103779
$ws.Cells.Item(5, 13).Value = 1.073413425707827
$ws.Cells.Item(5, 14).Value = 1.055458793998487
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.049633712286027
$ws.Cells.Item(6, 4).Value = 1.049987838496122
$ws.Cells.Item(6, 5).Value = 1.063333761254983
$ws.Cells.Item(6, 6).Value = 1.07111195520378
$ws.Cells.Item(6, 9).Value = 1.047928363479651
$ws.Cells.Item(6, 10).Value = 1.053989482485574
$ws.Cells.Item(6, 11).Value = 1.052377343209166
$ws.Cells.Item(6, 12).Value = 1.065691741615235
$ws.Cells.Item(6, 13).Value = 1.073451951118975
$ws.Cells.Item(6, 14).Value = 1.055486268444987
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.04933417814067
$ws.Cells.Item(7, 4).Value = 1.049763269283034
$ws.Cells.Item(7, 5).Value = 1.063036449426352
$ws.Cells.Item(7, 6).Value = 1.070788672475806
$ws.Cells.Item(7, 9).Value = 1.04783924192405
$ws.Cells.Item(7, 10).Value = 1.053800809454079
$ws.Cells.Item(7, 11).Value = 1.052212238143333
$ws.Cells.Item(7, 12).Value = 1.065453239352582
$ws.Cells.Item(7, 13).Value = 1.073187064406894
$ws.Cells.Item(7, 14).Value = 1.055297327476146
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.048082937771921
$ws.Cells.Item(8, 4).Value = 1.048825093597347
$ws.Cells.Item(8, 5).Value = 1.061795233371039
$ws.Cells.Item(8, 6).Value = 1.069439184916706
$ws.Cells.Item(8, 9).Value = 1.047464462461303
$ws.Cells.Item(8, 10).Value = 1.053011657944955
$ws.Cells.Item(8, 11).Value = 1.051521247370125
$ws.Cells.Item(8, 12).Value = 1.064456562979635
$ws.Cells.Item(8, 13).Value = 1.072080418235706
$ws.Cells.Item(8, 14).Value = 1.054507055281362
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.045881023972444
$ws.Cells.Item(9, 4).Value = 1.047173822183246
$ws.Cells.Item(9, 5).Value = 1.059613775050004
$ws.Cells.Item(9, 6).Value = 1.067068001420719
$ws.Cells.Item(9, 9).Value = 1.046795545326829
$ws.Cells.Item(9, 10).Value = 1.051619130578069
$ws.Cells.Item(9, 11).Value = 1.050300375236805
$ws.Cells.Item(9, 12).Value = 1.062701198124523
$ws.Cells.Item(9, 13).Value = 1.070132454315593
$ws.Cells.Item(9, 14).Value = 1.05311255036589
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.044415099474726
$ws.Cells.Item(10, 4).Value = 1.046074331394369
$ws.Cells.Item(10, 5).Value = 1.058163381996602
$ws.Cells.Item(10, 6).Value = 1.065491837441248
$ws.Cells.Item(10, 9).Value = 1.04634391016311
$ws.Cells.Item(10, 10).Value = 1.050689520218772
$ws.Cells.Item(10, 11).Value = 1.049484317549345
$ws.Cells.Item(10, 12).Value = 1.06153161850333
$ws.Cells.Item(10, 13).Value = 1.068835271849671
$ws.Cells.Item(10, 14).Value = 1.052181619853256
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.043780826036497
$ws.Cells.Item(11, 4).Value = 1.045598576063551
$ws.Cells.Item(11, 5).Value = 1.057536287043908
$ws.Cells.Item(11, 6).Value = 1.064810450769528
$ws.Cells.Item(11, 9).Value = 1.04614700353213
$ws.Cells.Item(11, 10).Value = 1.050286699629987
$ws.Cells.Item(11, 11).Value = 1.049130457067091
$ws.Cells.Item(11, 12).Value = 1.061025346112858
$ws.Cells.Item(11, 13).Value = 1.068273935325351
$ws.Cells.Item(11, 14).Value = 1.05177822721303
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.043545301680124
$ws.Cells.Item(12, 4).Value = 1.045421910283656
$ws.Cells.Item(12, 5).Value = 1.05730349752656
$ws.Cells.Item(12, 6).Value = 1.064557520070573
$ws.Cells.Item(12, 9).Value = 1.046073661925454
$ws.Cells.Item(12, 10).Value = 1.050137030806365
$ws.Cells.Item(12, 11).Value = 1.048998942839209
$ws.Cells.Item(12, 12).Value = 1.060837319387706
$ws.Cells.Item(12, 13).Value = 1.068065483712237
$ws.Cells.Item(12, 14).Value = 1.05162834584251
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.043595819137694
$ws.Cells.Item(13, 4).Value = 1.045459803370052
$ws.Cells.Item(13, 5).Value = 1.057353425249982
$ws.Cells.Item(13, 6).Value = 1.064611767007127
$ws.Cells.Item(13, 9).Value = 1.04608940306993
$ws.Cells.Item(13, 10).Value = 1.050169137225832
$ws.Cells.Item(13, 11).Value = 1.049027156455334
$ws.Cells.Item(13, 12).Value = 1.060877650596273
$ws.Cells.Item(13, 13).Value = 1.068110194843771
$ws.Cells.Item(13, 14).Value = 1.051660497856776
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.043761356025713
$ws.Cells.Item(14, 4).Value = 1.045583971769473
$ws.Cells.Item(14, 5).Value = 1.057517041697914
$ws.Cells.Item(14, 6).Value = 1.064789540033786
$ws.Cells.Item(14, 9).Value = 1.046140945204399
$ws.Cells.Item(14, 10).Value = 1.050274328837884
$ws.Cells.Item(14, 11).Value = 1.04911958757929
$ws.Cells.Item(14, 12).Value = 1.061009803242394
$ws.Cells.Item(14, 13).Value = 1.068256703548111
$ws.Cells.Item(14, 14).Value = 1.051765838852983
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.043863358465113
$ws.Cells.Item(15, 4).Value = 1.045660482801591
$ws.Cells.Item(15, 5).Value = 1.057617869968287
$ws.Cells.Item(15, 6).Value = 1.064899093974735
$ws.Cells.Item(15, 9).Value = 1.046172675298591
$ws.Cells.Item(15, 10).Value = 1.050339135143024
$ws.Cells.Item(15, 11).Value = 1.049176527570751
$ws.Cells.Item(15, 12).Value = 1.061091230238088
$ws.Cells.Item(15, 13).Value = 1.068346979546937
$ws.Cells.Item(15, 14).Value = 1.05183073719051
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.044457204318506
$ws.Cells.Item(16, 4).Value = 1.04610591278898
$ws.Cells.Item(16, 5).Value = 1.058205020006124
$ws.Cells.Item(16, 6).Value = 1.06553708214714
$ws.Cells.Item(16, 9).Value = 1.046356949895145
$ws.Cells.Item(16, 10).Value = 1.050716247979609
$ws.Cells.Item(16, 11).Value = 1.04950779158894
$ws.Cells.Item(16, 12).Value = 1.061565221647455
$ws.Cells.Item(16, 13).Value = 1.068872533397838
$ws.Cells.Item(16, 14).Value = 1.052208385570578
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.04482983754375
$ws.Cells.Item(17, 4).Value = 1.046385408790725
$ws.Cells.Item(17, 5).Value = 1.058573574710608
$ws.Cells.Item(17, 6).Value = 1.065937571193428
$ws.Cells.Item(17, 9).Value = 1.046472180534782
$ws.Cells.Item(17, 10).Value = 1.050952722919659
$ws.Cells.Item(17, 11).Value = 1.049715450852298
$ws.Cells.Item(17, 12).Value = 1.061862588254715
$ws.Cells.Item(17, 13).Value = 1.069202294312362
$ws.Cells.Item(17, 14).Value = 1.052445196332168
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.045047234614743
$ws.Cells.Item(18, 4).Value = 1.046548465982332
$ws.Cells.Item(18, 5).Value = 1.0587886365611
$ws.Cells.Item(18, 6).Value = 1.066171276010412
$ws.Cells.Item(18, 9).Value = 1.046539262753854
$ws.Cells.Item(18, 10).Value = 1.051090626416163
$ws.Cells.Item(18, 11).Value = 1.049836526489883
$ws.Cells.Item(18, 12).Value = 1.062036052865114
$ws.Cells.Item(18, 13).Value = 1.069394672182944
$ws.Cells.Item(18, 14).Value = 1.052583295667455
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.045121369249703
$ws.Cells.Item(19, 4).Value = 1.046604069647465
$ws.Cells.Item(19, 5).Value = 1.058861982350349
$ws.Cells.Item(19, 6).Value = 1.066250981319624
$ws.Cells.Item(19, 9).Value = 1.046562114016301
$ws.Cells.Item(19, 10).Value = 1.051137643116592
$ws.Cells.Item(19, 11).Value = 1.04987780190467
$ws.Cells.Item(19, 12).Value = 1.062095202451667
$ws.Cells.Item(19, 13).Value = 1.069460273790176
$ws.Cells.Item(19, 14).Value = 1.05263037913699
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.044789852719252
$ws.Cells.Item(20, 4).Value = 1.04635541820679
$ws.Cells.Item(20, 5).Value = 1.058534022949105
$ws.Cells.Item(20, 6).Value = 1.065894591503435
$ws.Cells.Item(20, 9).Value = 1.046459830806814
$ws.Cells.Item(20, 10).Value = 1.050927354322151
$ws.Cells.Item(20, 11).Value = 1.04969317598949
$ws.Cells.Item(20, 12).Value = 1.06183068200198
$ws.Cells.Item(20, 13).Value = 1.069166910573411
$ws.Cells.Item(20, 14).Value = 1.05241979170834
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.043712607516247
$ws.Cells.Item(21, 4).Value = 1.045547405855521
$ws.Cells.Item(21, 5).Value = 1.057468856818539
$ws.Cells.Item(21, 6).Value = 1.064737185700105
$ws.Cells.Item(21, 9).Value = 1.046125772895207
$ws.Cells.Item(21, 10).Value = 1.050243353720025
$ws.Cells.Item(21, 11).Value = 1.049092370972838
$ws.Cells.Item(21, 12).Value = 1.060970886874428
$ws.Cells.Item(21, 13).Value = 1.068213558895083
$ws.Cells.Item(21, 14).Value = 1.051734819746904
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.043035723632965
$ws.Cells.Item(22, 4).Value = 1.045039671666106
$ws.Cells.Item(22, 5).Value = 1.056799963603904
$ws.Cells.Item(22, 6).Value = 1.064010443403911
$ws.Cells.Item(22, 9).Value = 1.045914570580921
$ws.Cells.Item(22, 10).Value = 1.049813045109494
$ws.Cells.Item(22, 11).Value = 1.048714189256284
$ws.Cells.Item(22, 12).Value = 1.060430446829258
$ws.Cells.Item(22, 13).Value = 1.067614460362771
$ws.Cells.Item(22, 14).Value = 1.051303900048789
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.043394512264728
$ws.Cells.Item(23, 4).Value = 1.045308802800666
$ws.Cells.Item(23, 5).Value = 1.057154478577528
$ws.Cells.Item(23, 6).Value = 1.064395611478972
$ws.Cells.Item(23, 9).Value = 1.046026643377003
$ws.Cells.Item(23, 10).Value = 1.050041183327395
$ws.Cells.Item(23, 11).Value = 1.048914711214909
$ws.Cells.Item(23, 12).Value = 1.060716930078725
$ws.Cells.Item(23, 13).Value = 1.067932024009875
$ws.Cells.Item(23, 14).Value = 1.051532362249127
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.044807919996912
$ws.Cells.Item(24, 4).Value = 1.046368969562272
$ws.Cells.Item(24, 5).Value = 1.058551894409236
$ws.Cells.Item(24, 6).Value = 1.065914011846812
$ws.Cells.Item(24, 9).Value = 1.046465411518572
$ws.Cells.Item(24, 10).Value = 1.050938817386979
$ws.Cells.Item(24, 11).Value = 1.049703241191083
$ws.Cells.Item(24, 12).Value = 1.061845099015745
$ws.Cells.Item(24, 13).Value = 1.069182898857245
$ws.Cells.Item(24, 14).Value = 1.052431271052036
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.046449919233728
$ws.Cells.Item(25, 4).Value = 1.047600481394471
$ws.Cells.Item(25, 5).Value = 1.060177049456596
$ws.Cells.Item(25, 6).Value = 1.067680198587524
$ws.Cells.Item(25, 9).Value = 1.046969481602184
$ws.Cells.Item(25, 10).Value = 1.05197935757028
$ws.Cells.Item(25, 11).Value = 1.050616381145457
$ws.Cells.Item(25, 12).Value = 1.063154889032306
$ws.Cells.Item(25, 13).Value = 1.07063579645088
$ws.Cells.Item(25, 14).Value = 1.053473288921749

Write-Host "Applied vm_pu.xlsx updates for Case_4_215 (380 kV slack set to 1.02 pu)."
